$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("erosion")

# New rows of erosion data (continuing the existing series) for the "pure EM" run
$data = @(
    @(4,  2.789, 5),
    @(5,  2.514, 6),
    @(6,  2.298, 7),
    @(7,  2.103, 8),
    @(8,  1.987, 9),
    @(9,  1.766, 10),
    @(10, 1.625, 11),
    @(11, 1.489, 12),
    @(12, 1.379, 13),
    @(13, 1.32,  14),
    @(14, 1.256, 15),
    @(15, 1.201, 16),
    @(16, 1.155, 17),
    @(17, 1.106, 18),
    @(18, 0.979, 19),
    @(19, 0.897, 20)
)

$row = 6
foreach ($item in $data) {
    $ws1.Range("A$row").Value = $item[0]
    $ws1.Range("B$row").Value = $item[1]
    $ws1.Range("C$row").Value = $item[2]
    $ws1.Range("D$row").Value = "沈阳站"
    $row = $row + 1
}

# Make "erosion" the active sheet/tab with F9 selected
$ws1.Activate()
$ws1.Range("F9").Select()
